# Fix cx_freeze on mac
# Applies the leaderboard data corrections described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (was Teacher/Abhi/Not Applicable/Not Applicable -> Student/Aditya/7/E)
$ws.Range("A12").Value = "Student"
$ws.Range("B12").Value = "Aditya"
$ws.Range("C12").Value = "7"
$ws.Range("D12").Value = "E"

# Row 13 (only name changes: hy -> Abhi)
$ws.Range("B13").Value = "Abhi"

# Row 14 (was Student/Aditya/7/E -> Teacher/hy/Not Applicable/Not Applicable)
$ws.Range("A14").Value = "Teacher"
$ws.Range("B14").Value = "hy"
$ws.Range("C14").Value = "Not Applicable"
$ws.Range("D14").Value = "Not Applicable"

# Row 23 (only name changes: buck -> Abhijit)
$ws.Range("B23").Value = "Abhijit"

# Row 24 (only name changes: f -> buck)
$ws.Range("B24").Value = "buck"

# Row 25 (only name changes: Abhijit -> f)
$ws.Range("B25").Value = "f"

# Row 34 (only name changes: AA -> nk)
$ws.Range("B34").Value = "nk"

# Row 35 (only name changes: nk -> AA)
$ws.Range("B35").Value = "AA"

# Row 36 (was AAAA/N/A/N/A/4 -> mB/Not Applicable/Not Applicable/2)
$ws.Range("B36").Value = "mB"
$ws.Range("C36").Value = "Not Applicable"
$ws.Range("D36").Value = "Not Applicable"
$ws.Range("E36").Value = 2

# Row 37 (only name changes: mB -> Mama)
$ws.Range("B37").Value = "Mama"

# Row 38 (was Mama/Not Applicable/Not Applicable/2 -> A/N/A/N/A/2)
$ws.Range("B38").Value = "A"
$ws.Range("C38").Value = "N/A"
$ws.Range("D38").Value = "N/A"
